$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename "Testergebnis" header to "Ergebnis" while it is still row 1 (E1) ---
$ws.Range("E1").Value = "Ergebnis"

# --- Insert a new title row above the table ---
$ws.Rows("1").Insert()

# --- Title row: "Tests für HUD" merged across A1:E1 ---
$titleRng = $ws.Range("A1:E1")
$titleRng.Merge()
$titleRng.Value = "Tests für HUD"
$titleRng.Borders.LineStyle = 1
$titleRng.Font.Bold = $true
$titleRng.Font.Size = 16
$titleRng.HorizontalAlignment = -4108
$ws.Rows(1).RowHeight = 21

# --- Column widths (now narrower, fixed instead of best-fit) ---
$ws.Columns("B").ColumnWidth = 22.333333333333336
$ws.Columns("C").ColumnWidth = 19.166666666666668
$ws.Columns("D").ColumnWidth = 26.333333333333336
$ws.Columns("E").ColumnWidth = 10.666666666666666

# --- Row heights for the (now shifted) data rows ---
$ws.Rows(3).RowHeight = 93.6
$ws.Rows(4).RowHeight = 106.8
$ws.Rows(5).RowHeight = 101.4
$ws.Rows(6).RowHeight = 86.4

# --- Wrap text for the data cells (columns A-D), and vertical centering for A + E ---
$ws.Range("B3:D6").WrapText = $true
$ws.Range("A3:A6").WrapText = $true
$ws.Range("A3:A6").VerticalAlignment = -4108
$ws.Range("E3:E6").WrapText = $true
$ws.Range("E3:E6").VerticalAlignment = -4108

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
$ws.PageSetup.Zoom = 96
$ws.PageSetup.FitToPagesTall = $false

# --- Selection matches the saved state in the target file ---
$ws.Range("G4").Select()
